# Auto-generated edit script: updates cell values per the commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 2015.3529
$ws.Range("I80").Value = 334.4
$ws.Range("K80").Value = 1003.2
$ws.Range("M80").Value = -5.199999999999932

$ws.Range("H83").Value = 2015.3529
$ws.Range("I83").Value = 334.4
$ws.Range("K83").Value = 3009.6
$ws.Range("M83").Value = 1982.4

$ws.Range("H88").Value = 76928400
$ws.Range("I88").Value = 200003440
$ws.Range("J88").Value = 6499.75
$ws.Range("K88").Value = 200003440
$ws.Range("L88").Value = 6499.75
$ws.Range("M88").Value = -200003034
$ws.Range("N88").Value = -7311.75

$ws.Range("H91").Value = 76928400
$ws.Range("I91").Value = 200003440
$ws.Range("J91").Value = 6499.75
$ws.Range("K91").Value = 200003440
$ws.Range("L91").Value = 6499.75
$ws.Range("M91").Value = -200002036
$ws.Range("N91").Value = -9307.75

$ws.Range("H106").Value = 9499.143
$ws.Range("I106").Value = 5368.778
$ws.Range("K106").Value = 5368.778
$ws.Range("M106").Value = -4737.778

$ws.Range("H113").Value = 8056.875
$ws.Range("I113").Value = 4404.6665
$ws.Range("K113").Value = 4404.6665
$ws.Range("M113").Value = -1150.6665

$ws.Range("H138").Value = 3035.5615
$ws.Range("I138").Value = 1588.2593
$ws.Range("J138").Value = 4338.1333
$ws.Range("K138").Value = 4764.7779
$ws.Range("L138").Value = 13014.3999
$ws.Range("M138").Value = 375.2221
$ws.Range("N138").Value = -23294.3999

$ws.Range("H141").Value = 1223.5962
$ws.Range("I141").Value = 993.1087
$ws.Range("K141").Value = 2979.3261
$ws.Range("M141").Value = 2200.6739

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H46").Value = 49999
$ws.Range("J46").Value = 49999
$ws.Range("L46").Value = 49999
$ws.Range("N46").Value = -50637

$ws.Range("H61").Value = 6895.1577
$ws.Range("J61").Value = 10513.5
$ws.Range("L61").Value = 10513.5
$ws.Range("N61").Value = -10937.5

$ws.Range("H136").Value = 6895.1577
$ws.Range("J136").Value = 10513.5
$ws.Range("L136").Value = 31540.5
$ws.Range("N136").Value = -36640.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2350.2683
$ws.Range("J20").Value = 2903.85
$ws.Range("L20").Value = 2903.85
$ws.Range("N20").Value = -3397.85

$ws.Range("H86").Value = 3120.4443
$ws.Range("I86").Value = 1444.05
$ws.Range("K86").Value = 1444.05
$ws.Range("M86").Value = -321.05

$ws.Range("H89").Value = 3120.4443
$ws.Range("I89").Value = 1444.05
$ws.Range("K89").Value = 7220.25
$ws.Range("M89").Value = -1604.25

$ws.Range("H98").Value = 73541.664
$ws.Range("J98").Value = 73541.664
$ws.Range("L98").Value = 73541.664
$ws.Range("N98").Value = -79531.664

$ws.Range("H134").Value = 2267.7896
$ws.Range("I134").Value = 1598.375
$ws.Range("J134").Value = 5838
$ws.Range("K134").Value = 4795.125
$ws.Range("L134").Value = 17514
$ws.Range("M134").Value = -2260.125
$ws.Range("N134").Value = -22584

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2603.3076
$ws.Range("I16").Value = 2085.4285
$ws.Range("K16").Value = 2085.4285
$ws.Range("M16").Value = -1798.4285

$ws.Range("H31").Value = 16940.246
$ws.Range("I31").Value = 1953.0256
$ws.Range("J31").Value = 32321.87
$ws.Range("K31").Value = 1953.0256
$ws.Range("L31").Value = 32321.87
$ws.Range("M31").Value = -1658.0256
$ws.Range("N31").Value = -32911.87

$ws.Range("H34").Value = 16940.246
$ws.Range("I34").Value = 1953.0256
$ws.Range("J34").Value = 32321.87
$ws.Range("K34").Value = 1953.0256
$ws.Range("L34").Value = 32321.87
$ws.Range("M34").Value = -1751.0256
$ws.Range("N34").Value = -32725.87

$ws.Range("H62").Value = 14908.077
$ws.Range("J62").Value = 19444.889
$ws.Range("L62").Value = 19444.889
$ws.Range("N62").Value = -20692.889

$ws.Range("H65").Value = 14908.077
$ws.Range("J65").Value = 19444.889
$ws.Range("L65").Value = 97224.44499999999
$ws.Range("N65").Value = -103464.445

$ws.Range("H113").Value = 2603.3076
$ws.Range("I113").Value = 2085.4285
$ws.Range("K113").Value = 2085.4285
$ws.Range("M113").Value = 84.57150000000001

$ws.Range("H132").Value = 2672.1714
$ws.Range("I132").Value = 2318.5454
$ws.Range("K132").Value = 6955.6362
$ws.Range("M132").Value = -4425.6362

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H81").Value = 3358.3333
$ws.Range("J81").Value = 4062.4546
$ws.Range("L81").Value = 12187.3638
$ws.Range("N81").Value = -14433.3638

$ws.Range("H84").Value = 3358.3333
$ws.Range("J84").Value = 4062.4546
$ws.Range("L84").Value = 36562.0914
$ws.Range("N84").Value = -47794.0914

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 3176701.5
$ws.Range("I102").Value = 4763863.5
$ws.Range("J102").Value = 2377.8572
$ws.Range("K102").Value = 4763863.5
$ws.Range("L102").Value = 2377.8572
$ws.Range("M102").Value = -4762241.5
$ws.Range("N102").Value = -5621.8572

$ws.Range("H122").Value = 7653.3076
$ws.Range("I122").Value = 8338.950000000001
$ws.Range("J122").Value = 5367.8335
$ws.Range("K122").Value = 25016.85
$ws.Range("L122").Value = 16103.5005
$ws.Range("M122").Value = -22566.85
$ws.Range("N122").Value = -21003.5005

$ws.Range("H132").Value = 4497.69
$ws.Range("I132").Value = 3909.7334
$ws.Range("K132").Value = 11729.2002
$ws.Range("M132").Value = -9199.200199999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 4767.931
$ws.Range("I40").Value = 3930.4546
$ws.Range("K40").Value = 3930.4546
$ws.Range("M40").Value = -3794.4546

$ws.Range("H82").Value = 7772.8276
$ws.Range("I82").Value = 5233.636
$ws.Range("J82").Value = 9324.556
$ws.Range("K82").Value = 5233.636
$ws.Range("L82").Value = 9324.556
$ws.Range("M82").Value = -4872.636
$ws.Range("N82").Value = -10046.556

$ws.Range("H85").Value = 7772.8276
$ws.Range("I85").Value = 5233.636
$ws.Range("J85").Value = 9324.556
$ws.Range("K85").Value = 5233.636
$ws.Range("L85").Value = 9324.556
$ws.Range("M85").Value = -3985.636
$ws.Range("N85").Value = -11820.556

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 7500.25
$ws.Range("I62").Value = 7500.25
$ws.Range("K62").Value = 7500.25
$ws.Range("M62").Value = -6876.25

$ws.Range("H65").Value = 7500.25
$ws.Range("I65").Value = 7500.25
$ws.Range("K65").Value = 37501.25
$ws.Range("M65").Value = -34381.25

$ws.Range("H100").Value = 294
$ws.Range("I100").Value = 294
$ws.Range("K100").Value = 588
$ws.Range("M100").Value = -47

$ws.Range("H132").Value = 2331.6584
$ws.Range("I132").Value = 1575.4412
$ws.Range("J132").Value = 6004.7144
$ws.Range("K132").Value = 4726.3236
$ws.Range("L132").Value = 18014.1432
$ws.Range("M132").Value = -2196.3236
$ws.Range("N132").Value = -23074.1432
